# Insert a new "Date"-styled paragraph ("Version 0.1 (25.06.2026)") right
# after the Subtitle paragraph ("KI-unterstützte Dokumentation") and before
# the Table-of-Contents content control (sdt) that follows it.
#
# Quirk: inserting a paragraph break exactly at the boundary between the
# Subtitle paragraph and the following Table-of-Contents sdt always lands
# the new paragraph *inside* the sdt's content instead of as a body-level
# sibling before it - true for InsertParagraphBefore/After and for InsertXML
# alike. To work around this we create the new paragraph break at the
# *other* (safe, ordinary) side of the Subtitle paragraph - i.e. just
# before it - which gives us two paragraphs in the right place count-wise,
# then fill them with the right content/styles via InsertXML (which lets us
# supply the <w:pPr> directly, so no stray rsid/formatting artifacts get
# stamped) and delete the now-redundant original paragraph.

$d = $word.ActiveDocument

$subtitle = $d.Paragraphs.Item(2)
if ($subtitle.Range.Text.TrimEnd([char]13) -ne "KI-unterstützte Dokumentation") {
    throw "unexpected paragraph 2 content: $($subtitle.Range.Text)"
}

# Make room: insert an empty paragraph (inherits the Subtitle style) right
# before the Subtitle paragraph. This boundary is an ordinary one (not
# adjacent to the sdt) so it behaves normally.
$insertionPoint = $subtitle.Range.Duplicate
$insertionPoint.Collapse(1)
$insertionPoint.InsertParagraphBefore()

# Paragraph 2 is now the new empty paragraph; paragraph 3 is the original,
# untouched "KI-unterstützte Dokumentation" paragraph.
$newPara = $d.Paragraphs.Item(2)

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr><w:r><w:t xml:space="preserve">KI-unterstützte Dokumentation</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="Date"/></w:pPr><w:r><w:t xml:space="preserve">Version 0.1 (25.06.2026)</w:t></w:r></w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

# Fill the new (empty) paragraph's range with both real paragraphs at once -
# this replaces the single empty paragraph with the two fully-formed ones.
$null = $newPara.Range.InsertXML($packageXml)

# The original "KI-unterstützte Dokumentation" paragraph (now shifted down
# to index 4) is a now-redundant duplicate - remove it.
$dup = $d.Paragraphs.Item(4)
if ($dup.Range.Text.TrimEnd([char]13) -ne "KI-unterstützte Dokumentation") {
    throw "unexpected duplicate paragraph content: $($dup.Range.Text)"
}
$dup.Range.Delete()
